# Apply the timesheet update to the "29-04-2022" sheet (the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 / Row 8: lunch break boundary moved from 13:45 to 13:40 ---
$ws.Range("E7").Value = 0.56944444444444442
$ws.Range("D8").Value = 0.56944444444444442

# --- Row 72: task name relabelled ---
$ws.Range("B72").Value = "Worked on HTML layout(Achivements)"

# --- Rows 77-89: afternoon/evening task block rewritten ---
$ws.Range("B77").Value = "Modifying MyProfile page"
$ws.Range("C77").Value = "Project"
$ws.Range("D77").Value = 0.36458333333333331
$ws.Range("E77").Value = 0.39583333333333331

$ws.Range("B78").Value = "Soft Skill"
$ws.Range("C78").Value = "Non Project"
$ws.Range("D78").Value = 0.39583333333333331
$ws.Range("E78").Value = 0.4375

$ws.Range("B79").Value = "Morning Break"
$ws.Range("C79").Value = "Lunch and Break"
$ws.Range("D79").Value = 0.4375
$ws.Range("E79").Value = 0.45833333333333331

$ws.Range("B80").Value = "Modifying MyProfile page"
$ws.Range("C80").Value = "Project"
$ws.Range("D80").Value = 0.45833333333333331
$ws.Range("E80").Value = 0.5

$ws.Range("B81").Value = "Timesheet Entering"
$ws.Range("C81").Value = "Non Project"
$ws.Range("D81").Value = 0.5
$ws.Range("E81").Value = 0.52083333333333337

$ws.Range("B82").Value = "Customer meeting"
$ws.Range("C82").Value = "Project"
$ws.Range("D82").Value = 0.52083333333333337
$ws.Range("E82").Value = 0.55555555555555558

$ws.Range("B83").Value = "Modifying MyProfile page"
$ws.Range("C83").Value = "Project"
$ws.Range("D83").Value = 0.58333333333333337
$ws.Range("E83").Value = 0.625

$ws.Range("B84").Value = "Modifying create Page"
$ws.Range("C84").Value = "Project"
$ws.Range("D84").Value = 0.625
$ws.Range("E84").Value = 0.65625

$ws.Range("B85").Value = "Break"
$ws.Range("C85").Value = "Lunch and Break"
$ws.Range("D85").Value = 0.66666666666666663
$ws.Range("E85").Value = 0.6875

$ws.Range("B86").Value = "Modifying create Page"
$ws.Range("C86").Value = "Project"
$ws.Range("D86").Value = 0.6875
$ws.Range("E86").Value = 0.70833333333333337

$ws.Range("B87").Value = "Team Meeting"
$ws.Range("C87").Value = "Meeting "
$ws.Range("D87").Value = 0.70833333333333337
$ws.Range("E87").Value = 0.73958333333333337

$ws.Range("B88").Value = "Modifying MyProfile page"
$ws.Range("C88").Value = "Project"
$ws.Range("D88").Value = 0.89583333333333337
$ws.Range("E88").Value = 0.95833333333333337
$ws.Range("F88").Formula = "=E88-D88"

$ws.Range("B89").Value = "Exploration on web API"
$ws.Range("C89").Value = "Non Project"
$ws.Range("D89").Value = 0.95833333333333337
$ws.Range("E89").Value = 0.99305555555555547
$ws.Range("F89").Formula = "=E89-D89"

# --- Sheet view: active selection moved ---
$ws.Range("F89").Select()
